$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right marking 5 -> 4, Wrong marking -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 120 -> 96, Wrong total -2 -> -4
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -4

# E12: "120 / 140" -> "92 / 112"
$ws.Range("E12").Value = "92 / 112"
